$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H52").Value = 1500
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 1500
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 4500
$ws.Range("N52").Value = -4820
$ws.Range("M52").ClearContents()

$ws.Range("H125").Value = 3923495
$ws.Range("I125").Value = 2890
$ws.Range("J125").Value = 11764705
$ws.Range("K125").Value = 26010
$ws.Range("L125").Value = 105882345
$ws.Range("M125").Value = -23550
$ws.Range("N125").Value = -105887265

$ws.Range("H141").Value = 8832.5
$ws.Range("I141").Value = 1413.4615
$ws.Range("J141").Value = 17600.455
$ws.Range("K141").Value = 4240.3845
$ws.Range("L141").Value = 52801.36500000001
$ws.Range("M141").Value = 939.6154999999999
$ws.Range("N141").Value = -63161.36500000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 26338.334
$ws.Range("I21").Value = 15
$ws.Range("J21").Value = 39500
$ws.Range("K21").Value = 15
$ws.Range("L21").Value = 39500
$ws.Range("M21").Value = 359
$ws.Range("N21").Value = -40248

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2631.125
$ws.Range("I20").Value = 2577.8333
$ws.Range("K20").Value = 2577.8333
$ws.Range("M20").Value = -2330.8333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1167.3704
$ws.Range("I16").Value = 1214.45
$ws.Range("J16").Value = 1032.8572
$ws.Range("K16").Value = 1214.45
$ws.Range("L16").Value = 1032.8572
$ws.Range("M16").Value = -927.45
$ws.Range("N16").Value = -1606.8572

$ws.Range("H107").Value = 499.8889
$ws.Range("I107").Value = 280.92856
$ws.Range("J107").Value = 735.6923
$ws.Range("K107").Value = 280.92856
$ws.Range("L107").Value = 735.6923
$ws.Range("M107").Value = 1639.07144
$ws.Range("N107").Value = -4575.6923

$ws.Range("H113").Value = 1167.3704
$ws.Range("I113").Value = 1214.45
$ws.Range("J113").Value = 1032.8572
$ws.Range("K113").Value = 1214.45
$ws.Range("L113").Value = 1032.8572
$ws.Range("M113").Value = 955.55
$ws.Range("N113").Value = -5372.8572

$ws.Range("H132").Value = 2741.5925
$ws.Range("I132").Value = 1995.6666
$ws.Range("J132").Value = 3674
$ws.Range("K132").Value = 5986.9998
$ws.Range("L132").Value = 11022
$ws.Range("M132").Value = -3456.9998
$ws.Range("N132").Value = -16082

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 59.38095
$ws.Range("I12").Value = 7
$ws.Range("K12").Value = 21
$ws.Range("M12").Value = 152

$ws.Range("H98").Value = 1100
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 1100
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 3300
$ws.Range("N98").Value = -6296
$ws.Range("M98").ClearContents()

$ws.Range("H108").Value = 6800
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 6800
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 20400
$ws.Range("N108").Value = -26160
$ws.Range("M108").ClearContents()

$ws.Range("H110").Value = 5000
$ws.Range("I110").Value = 1000
$ws.Range("J110").Value = 7000
$ws.Range("K110").Value = 3000
$ws.Range("L110").Value = 21000
$ws.Range("M110").Value = 1090
$ws.Range("N110").Value = -29180

$ws.Range("H113").Value = 1081.6666
$ws.Range("I113").Value = 558
$ws.Range("J113").Value = 1500.6
$ws.Range("K113").Value = 1674
$ws.Range("L113").Value = 4501.799999999999
$ws.Range("M113").Value = 496
$ws.Range("N113").Value = -8841.799999999999

$ws.Range("H114").Value = 392.375
$ws.Range("I114").Value = 162.71428
$ws.Range("J114").Value = 2000
$ws.Range("K114").Value = 488.14284
$ws.Range("L114").Value = 6000
$ws.Range("M114").Value = 2765.85716
$ws.Range("N114").Value = -12508

$ws.Range("H117").Value = 1750
$ws.Range("J117").Value = 3000
$ws.Range("L117").Value = 9000
$ws.Range("N117").Value = -15884

$ws.Range("H120").Value = 3758.25
$ws.Range("J120").Value = 10033
$ws.Range("L120").Value = 30099
$ws.Range("N120").Value = -39775

$ws.Range("H123").Value = 2570.4285
$ws.Range("I123").Value = 780
$ws.Range("J123").Value = 3286.6
$ws.Range("K123").Value = 2340
$ws.Range("L123").Value = 9859.799999999999
$ws.Range("M123").Value = 110
$ws.Range("N123").Value = -14759.8

$ws.Range("H124").Value = 1065
$ws.Range("I124").Value = 530
$ws.Range("K124").Value = 1590
$ws.Range("M124").Value = 3320

$ws.Range("H131").Value = 1105.6044
$ws.Range("I131").Value = 765.4
$ws.Range("J131").Value = 1125.3837
$ws.Range("K131").Value = 2296.2
$ws.Range("L131").Value = 3376.1511
$ws.Range("M131").Value = 2743.8
$ws.Range("N131").Value = -13456.1511

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H124").Value = 38500
$ws.Range("J124").Value = 38500
$ws.Range("L124").Value = 38500
$ws.Range("N124").Value = -48320

$ws.Range("H126").Value = 1953.3143
$ws.Range("I126").Value = 1615
$ws.Range("J126").Value = 2799.1
$ws.Range("K126").Value = 4845
$ws.Range("L126").Value = 8397.299999999999
$ws.Range("M126").Value = -2375
$ws.Range("N126").Value = -13337.3

$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

$ws.Range("H132").Value = 4109.35
$ws.Range("I132").Value = 2751.923
$ws.Range("J132").Value = 6630.2856
$ws.Range("K132").Value = 8255.769
$ws.Range("L132").Value = 19890.8568
$ws.Range("M132").Value = -5725.769
$ws.Range("N132").Value = -24950.8568

$ws.Range("H140").Value = 41500
$ws.Range("J140").Value = 41500
$ws.Range("L140").Value = 41500
$ws.Range("N140").Value = -51860

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1989.5
$ws.Range("I68").Value = 1500
$ws.Range("J68").Value = 2087.4
$ws.Range("K68").Value = 1500
$ws.Range("L68").Value = 2087.4
$ws.Range("M68").Value = -751
$ws.Range("N68").Value = -3585.4

$ws.Range("H71").Value = 1989.5
$ws.Range("I71").Value = 1500
$ws.Range("J71").Value = 2087.4
$ws.Range("K71").Value = 7500
$ws.Range("L71").Value = 10437
$ws.Range("M71").Value = -3756
$ws.Range("N71").Value = -17925

$ws.Range("H119").Value = 26613.334
$ws.Range("J119").Value = 26613.334
$ws.Range("L119").Value = 26613.334
$ws.Range("N119").Value = -36289.334

$ws.Range("H132").Value = 5110.5
$ws.Range("I132").Value = 4298.2
$ws.Range("J132").Value = 6125.875
$ws.Range("K132").Value = 12894.6
$ws.Range("L132").Value = 18377.625
$ws.Range("M132").Value = -10364.6
$ws.Range("N132").Value = -23437.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2993.3333
$ws.Range("I62").Value = 2592
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 2592
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -1968
$ws.Range("N62").Value = -6248

$ws.Range("H65").Value = 2993.3333
$ws.Range("I65").Value = 2592
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 12960
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -9840
$ws.Range("N65").Value = -31240

$ws.Range("H107").Value = 1567.4348
$ws.Range("I107").Value = 1256.35
$ws.Range("J107").Value = 3641.3333
$ws.Range("K107").Value = 3769.05
$ws.Range("L107").Value = 10923.9999
$ws.Range("M107").Value = -1849.05
$ws.Range("N107").Value = -14763.9999

$ws.Range("H132").Value = 2741.2942
$ws.Range("I132").Value = 2655.4443
$ws.Range("J132").Value = 2837.875
$ws.Range("K132").Value = 7966.3329
$ws.Range("L132").Value = 8513.625
$ws.Range("M132").Value = -5436.3329
$ws.Range("N132").Value = -13573.625
